$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 currently holds the shared string "R40" (t="s", s="23").
# The target edit changes it to hold the literal text "1" -- as a real
# text/shared-string value (not a number) -- while leaving its style (s="23")
# and every other cell/format completely untouched.
#
# A direct  $ws.Range("B11").Value = "1"  assignment gets auto-coerced to the
# *number* 1 by the engine (just like typing 1 into a General formatted cell
# in real Excel would), and forcing text via NumberFormat="@" or a leading
# apostrophe both mint a brand new cell style, which would incorrectly change
# B11's "s" attribute. To avoid both problems we stage the text value in a
# scratch cell far outside the used range, build it there as a text formula
# result, then copy/paste-special just the *value* (xlPasteValues = -4163)
# onto B11. PasteSpecial values-only preserves the destination's existing
# style while bringing over the source's text typing, which is exactly the
# semantics we need. Finally the scratch cell is fully cleared (contents +
# formatting) so it leaves no trace in the saved workbook.

$helper = $ws.Range("Z1")
$helper.Formula = '="1"'
$helper.Copy()

$target = $ws.Range("B11")
$target.PasteSpecial(-4163)

$helper.ClearContents()
$helper.Clear()
